# Apply text corrections to the RF018 - Gerenciar Metas de Desempenho test case sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Precondition text: remove the "; e," and replace with " e" (both occurrences, TC1 and TC2 blocks)
$oldPrecondition = "Lider de Pessoas esta autenticado no sistema; e, tem permissao para gerenciar Metas de Desempenho"
$newPrecondition = "Lider de Pessoas esta autenticado no sistema e tem permissao para gerenciar Metas de Desempenho"

# 2) Step text: add "para " before "modificar"
$oldStep = "Lider de Pessoas com uma avaliacao selecionada, clica na opcao 'Editar' modificar a Avaliacao de Desempenho"
$newStep = "Lider de Pessoas com uma avaliacao selecionada, clica na opcao 'Editar' para modificar a Avaliacao de Desempenho"

# 3) Expected result text: "constando" -> "contendo"
$oldExpected = "SYSTEM apresenta o formulario com o campo 'Metas' constando cada Competencia do perfil avaliado"
$newExpected = "SYSTEM apresenta o formulario com o campo 'Metas' contendo cada Competencia do perfil avaliado"

$used = $ws.UsedRange
foreach ($cell in $used.Cells) {
    $val = $cell.Value2
    if ($val -eq $oldPrecondition) {
        $cell.Value2 = $newPrecondition
    } elseif ($val -eq $oldStep) {
        $cell.Value2 = $newStep
    } elseif ($val -eq $oldExpected) {
        $cell.Value2 = $newExpected
    }
}
